$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 173 (pushes the existing row 173 and everything
# below it down by one, e.g. old row 173 -> new row 174, ..., old row 284 -> new row 285)
$ws.Rows.Item(173).Insert()

# Populate the newly inserted row 173 with a new "Cilantro" price record
$ws.Range("A173").Value = 3
$ws.Range("B173").Value = "Femacal de La Calera"
$ws.Range("C173").Value = "Coquimbo"
$ws.Range("D173").Value = 44582
$ws.Range("E173").Value = 5
$ws.Range("F173").Value = 100112040
$ws.Range("G173").Value = "Cilantro"
$ws.Range("H173").Value = "Sin especificar"
$ws.Range("I173").Value = "Primera"
$ws.Range("J173").Value = 120
$ws.Range("K173").Value = 4000
$ws.Range("L173").Value = 4500
$ws.Range("M173").Value = 4250
$ws.Range("N173").Value = "`$/docena de atados (3 kilos)"
$ws.Range("O173").Value = "Provincia de Quillota"
$ws.Range("P173").Value = 1417
$ws.Range("Q173").Value = 3
$ws.Range("R173").Value = "Hortaliza"
